$wb = $excel.ActiveWorkbook

$oldVersion = "mines - January 30 (built on February 02 2026 12.49.33 EST)"
$newVersion = "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"

# --- "About" sheet ---
$wsAbout = $wb.Worksheets.Item("About")

$a2 = $wsAbout.Range("A2")
$a2Text = $a2.Value()
$a2.Value = $a2Text.Replace($oldVersion, $newVersion)

$a6 = $wsAbout.Range("A6")
$a6Text = $a6.Value()
$a6.Value = $a6Text.Replace($oldVersion, $newVersion)

# --- "Boundaries and methane sources" sheet ---
$wsData = $wb.Worksheets.Item("Boundaries and methane sources")

$usedRange = $wsData.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $wsData.Cells.Item($r, 19)  # column S = 19
    $cellText = $cell.Value()
    if ($cellText -ne $null -and $cellText -ne "") {
        $cell.Value = $cellText.Replace($oldVersion, $newVersion)
    }
}
